$wb = $excel.ActiveWorkbook

# --- Step 1: set up sheet order / names -----------------------------------
# The current last sheet ("总计") is repurposed to hold the new quarter's
# fund table (it keeps its original sheetId) and a brand-new "总计" sheet is
# appended after it to hold the refreshed summary table.
$oldTotal = $wb.Worksheets.Item("总计")
$newTotal = $wb.Worksheets.Add($null, $oldTotal)

# Recreate the look of a "normal" sheet (outline + page margins) on the
# freshly inserted sheet to match the rest of the workbook.
$newTotal.Outline.SummaryRow = 1
$newTotal.Outline.SummaryColumn = 1
$newTotal.PageSetup.LeftMargin = 54
$newTotal.PageSetup.RightMargin = 54
$newTotal.PageSetup.TopMargin = 72
$newTotal.PageSetup.BottomMargin = 72
$newTotal.PageSetup.HeaderMargin = 36
$newTotal.PageSetup.FooterMargin = 36

# Grab the "总计" summary header (日期 / 持有数量(只) / 持有市值(亿元)) and the
# styled A-column look from the sheet before it gets overwritten below.
$oldTotal.Range("B1:D1").Copy($newTotal.Range("B1:D1"))
$oldTotal.Range("A2").Copy()
$newTotal.Range("A2:A6").PasteSpecial(-4122)

$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

$q1 = $wb.Worksheets.Item("2022-Q1")
$total = $wb.Worksheets.Item("总计")

# --- Step 2: header row + fund table on the "2022-Q1" sheet ---------------
# Re-use the already-styled header row from "2021-Q4" (same 7 columns) so the
# fonts/borders/alignment match the other quarterly sheets exactly.
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("B1:H1").Copy($q1.Range("B1:H1"))
$q4.Range("A2").Copy()
$q1.Range("A2:A13").PasteSpecial(-4122)

# row 2: 景顺长城中证500指数增强
$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).NumberFormat = "@"
$q1.Cells.Item(2, 2).Value = "006682"
$q1.Cells.Item(2, 3).Value = "景顺长城中证500指数增强"
$q1.Cells.Item(2, 4).NumberFormat = "@"
$q1.Cells.Item(2, 4).Value = "16.63"
$q1.Cells.Item(2, 5).NumberFormat = "@"
$q1.Cells.Item(2, 5).Value = "87.75"
$q1.Cells.Item(2, 6).NumberFormat = "@"
$q1.Cells.Item(2, 6).Value = "1.98"
$q1.Cells.Item(2, 7).NumberFormat = "@"
$q1.Cells.Item(2, 7).Value = "0.3293"
$q1.Cells.Item(2, 8).Value = 7

# row 3: 浙商中证500指数增强A
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).NumberFormat = "@"
$q1.Cells.Item(3, 2).Value = "002076"
$q1.Cells.Item(3, 3).Value = "浙商中证500指数增强A"
$q1.Cells.Item(3, 4).NumberFormat = "@"
$q1.Cells.Item(3, 4).Value = "14.53"
$q1.Cells.Item(3, 5).NumberFormat = "@"
$q1.Cells.Item(3, 5).Value = "93.68"
$q1.Cells.Item(3, 6).NumberFormat = "@"
$q1.Cells.Item(3, 6).Value = "1.47"
$q1.Cells.Item(3, 7).NumberFormat = "@"
$q1.Cells.Item(3, 7).Value = "0.2136"
$q1.Cells.Item(3, 8).Value = 6

# row 4: 景顺长城量化精选股票
$q1.Cells.Item(4, 1).Value = 2
$q1.Cells.Item(4, 2).NumberFormat = "@"
$q1.Cells.Item(4, 2).Value = "000978"
$q1.Cells.Item(4, 3).Value = "景顺长城量化精选股票"
$q1.Cells.Item(4, 4).NumberFormat = "@"
$q1.Cells.Item(4, 4).Value = "8.51"
$q1.Cells.Item(4, 5).NumberFormat = "@"
$q1.Cells.Item(4, 5).Value = "93.86"
$q1.Cells.Item(4, 6).NumberFormat = "@"
$q1.Cells.Item(4, 6).Value = "2.00"
$q1.Cells.Item(4, 7).NumberFormat = "@"
$q1.Cells.Item(4, 7).Value = "0.1702"
$q1.Cells.Item(4, 8).Value = 3

# row 5: 国投瑞银中证500指数量化增强A
$q1.Cells.Item(5, 1).Value = 3
$q1.Cells.Item(5, 2).NumberFormat = "@"
$q1.Cells.Item(5, 2).Value = "005994"
$q1.Cells.Item(5, 3).Value = "国投瑞银中证500指数量化增强A"
$q1.Cells.Item(5, 4).NumberFormat = "@"
$q1.Cells.Item(5, 4).Value = "11.53"
$q1.Cells.Item(5, 5).NumberFormat = "@"
$q1.Cells.Item(5, 5).Value = "87.00"
$q1.Cells.Item(5, 6).NumberFormat = "@"
$q1.Cells.Item(5, 6).Value = "1.46"
$q1.Cells.Item(5, 7).NumberFormat = "@"
$q1.Cells.Item(5, 7).Value = "0.1683"
$q1.Cells.Item(5, 8).Value = 5

# row 6: 景顺长城量化对冲策略三个月定期开放灵活配置混合
$q1.Cells.Item(6, 1).Value = 4
$q1.Cells.Item(6, 2).NumberFormat = "@"
$q1.Cells.Item(6, 2).Value = "008851"
$q1.Cells.Item(6, 3).Value = "景顺长城量化对冲策略三个月定期开放灵活配置混合"
$q1.Cells.Item(6, 4).NumberFormat = "@"
$q1.Cells.Item(6, 4).Value = "5.05"
$q1.Cells.Item(6, 5).NumberFormat = "@"
$q1.Cells.Item(6, 5).Value = "74.55"
$q1.Cells.Item(6, 6).NumberFormat = "@"
$q1.Cells.Item(6, 6).Value = "1.58"
$q1.Cells.Item(6, 7).NumberFormat = "@"
$q1.Cells.Item(6, 7).Value = "0.0798"
$q1.Cells.Item(6, 8).Value = 6

# row 7: 申万菱信量化小盘股票(LOF)
$q1.Cells.Item(7, 1).Value = 5
$q1.Cells.Item(7, 2).NumberFormat = "@"
$q1.Cells.Item(7, 2).Value = "163110"
$q1.Cells.Item(7, 3).Value = "申万菱信量化小盘股票(LOF)"
$q1.Cells.Item(7, 4).NumberFormat = "@"
$q1.Cells.Item(7, 4).Value = "5.68"
$q1.Cells.Item(7, 5).NumberFormat = "@"
$q1.Cells.Item(7, 5).Value = "92.25"
$q1.Cells.Item(7, 6).NumberFormat = "@"
$q1.Cells.Item(7, 6).Value = "1.35"
$q1.Cells.Item(7, 7).NumberFormat = "@"
$q1.Cells.Item(7, 7).Value = "0.0767"
$q1.Cells.Item(7, 8).Value = 2

# row 8: 国投瑞银中证500指数量化增强C
$q1.Cells.Item(8, 1).Value = 6
$q1.Cells.Item(8, 2).NumberFormat = "@"
$q1.Cells.Item(8, 2).Value = "007089"
$q1.Cells.Item(8, 3).Value = "国投瑞银中证500指数量化增强C"
$q1.Cells.Item(8, 4).NumberFormat = "@"
$q1.Cells.Item(8, 4).Value = "3.82"
$q1.Cells.Item(8, 5).NumberFormat = "@"
$q1.Cells.Item(8, 5).Value = "87.00"
$q1.Cells.Item(8, 6).NumberFormat = "@"
$q1.Cells.Item(8, 6).Value = "1.46"
$q1.Cells.Item(8, 7).NumberFormat = "@"
$q1.Cells.Item(8, 7).Value = "0.0558"
$q1.Cells.Item(8, 8).Value = 5

# row 9: 浙商中证500指数增强C
$q1.Cells.Item(9, 1).Value = 7
$q1.Cells.Item(9, 2).NumberFormat = "@"
$q1.Cells.Item(9, 2).Value = "007386"
$q1.Cells.Item(9, 3).Value = "浙商中证500指数增强C"
$q1.Cells.Item(9, 4).NumberFormat = "@"
$q1.Cells.Item(9, 4).Value = "3.38"
$q1.Cells.Item(9, 5).NumberFormat = "@"
$q1.Cells.Item(9, 5).Value = "93.68"
$q1.Cells.Item(9, 6).NumberFormat = "@"
$q1.Cells.Item(9, 6).Value = "1.47"
$q1.Cells.Item(9, 7).NumberFormat = "@"
$q1.Cells.Item(9, 7).Value = "0.0497"
$q1.Cells.Item(9, 8).Value = 6

# row 10: 景顺长城量化平衡灵活配置混合
$q1.Cells.Item(10, 1).Value = 8
$q1.Cells.Item(10, 2).NumberFormat = "@"
$q1.Cells.Item(10, 2).Value = "005258"
$q1.Cells.Item(10, 3).Value = "景顺长城量化平衡灵活配置混合"
$q1.Cells.Item(10, 4).NumberFormat = "@"
$q1.Cells.Item(10, 4).Value = "2.39"
$q1.Cells.Item(10, 5).NumberFormat = "@"
$q1.Cells.Item(10, 5).Value = "90.00"
$q1.Cells.Item(10, 6).NumberFormat = "@"
$q1.Cells.Item(10, 6).Value = "1.66"
$q1.Cells.Item(10, 7).NumberFormat = "@"
$q1.Cells.Item(10, 7).Value = "0.0397"
$q1.Cells.Item(10, 8).Value = 7

# row 11: 国投瑞银安睿混合A
$q1.Cells.Item(11, 1).Value = 9
$q1.Cells.Item(11, 2).NumberFormat = "@"
$q1.Cells.Item(11, 2).Value = "011731"
$q1.Cells.Item(11, 3).Value = "国投瑞银安睿混合A"
$q1.Cells.Item(11, 4).NumberFormat = "@"
$q1.Cells.Item(11, 4).Value = "2.58"
$q1.Cells.Item(11, 5).NumberFormat = "@"
$q1.Cells.Item(11, 5).Value = "43.48"
$q1.Cells.Item(11, 6).NumberFormat = "@"
$q1.Cells.Item(11, 6).Value = "0.72"
$q1.Cells.Item(11, 7).NumberFormat = "@"
$q1.Cells.Item(11, 7).Value = "0.0186"
$q1.Cells.Item(11, 8).Value = 6

# row 12: 国投瑞银安睿混合C
$q1.Cells.Item(12, 1).Value = 10
$q1.Cells.Item(12, 2).NumberFormat = "@"
$q1.Cells.Item(12, 2).Value = "011732"
$q1.Cells.Item(12, 3).Value = "国投瑞银安睿混合C"
$q1.Cells.Item(12, 4).NumberFormat = "@"
$q1.Cells.Item(12, 4).Value = "0.95"
$q1.Cells.Item(12, 5).NumberFormat = "@"
$q1.Cells.Item(12, 5).Value = "43.48"
$q1.Cells.Item(12, 6).NumberFormat = "@"
$q1.Cells.Item(12, 6).Value = "0.72"
$q1.Cells.Item(12, 7).NumberFormat = "@"
$q1.Cells.Item(12, 7).Value = "0.0068"
$q1.Cells.Item(12, 8).Value = 6

# row 13: 景顺长城量化先锋混合
$q1.Cells.Item(13, 1).Value = 11
$q1.Cells.Item(13, 2).NumberFormat = "@"
$q1.Cells.Item(13, 2).Value = "006201"
$q1.Cells.Item(13, 3).Value = "景顺长城量化先锋混合"
$q1.Cells.Item(13, 4).NumberFormat = "@"
$q1.Cells.Item(13, 4).Value = "0.09"
$q1.Cells.Item(13, 5).NumberFormat = "@"
$q1.Cells.Item(13, 5).Value = "46.30"
$q1.Cells.Item(13, 6).NumberFormat = "@"
$q1.Cells.Item(13, 6).Value = "0.85"
$q1.Cells.Item(13, 7).NumberFormat = "@"
$q1.Cells.Item(13, 7).Value = "0.0008"
$q1.Cells.Item(13, 8).Value = 8

# --- Step 3: refreshed summary rows on the new "总计" sheet ----------------
# row 2: 2022-Q1
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 12
$total.Cells.Item(2, 4).Value = 1.21

# row 3: 2021-Q4
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q4"
$total.Cells.Item(3, 3).Value = 16
$total.Cells.Item(3, 4).Value = 1.37

# row 4: 2021-Q3
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q3"
$total.Cells.Item(4, 3).Value = 9
$total.Cells.Item(4, 4).Value = 0.93

# row 5: 2021-Q2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(5, 2).Value = "2021-Q2"
$total.Cells.Item(5, 3).Value = 1
$total.Cells.Item(5, 4).Value = 0.02

# row 6: 2021-Q1
$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(6, 2).Value = "2021-Q1"
$total.Cells.Item(6, 3).Value = 7
$total.Cells.Item(6, 4).Value = 0.09

